$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H - copy formatting from neighboring header cell (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# "Save" indicator values per row (1 = saved, 0 = not saved)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
